$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3509
$ws.Range("L3").Value = 3669
$ws.Range("L4").Value = 912
$ws.Range("L5").Value = 220
$ws.Range("L6").Value = 3209
$ws.Range("L7").Value = 11519

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 389
$ws.Range("L8").Value = 747
$ws.Range("L9").Value = 73
$ws.Range("L10").Value = 72
$ws.Range("L11").Value = 192
$ws.Range("L13").Value = 18
$ws.Range("L15").Value = 83
$ws.Range("L17").Value = 21
$ws.Range("L18").Value = 89
$ws.Range("L19").Value = 326
$ws.Range("L20").Value = 291
$ws.Range("L23").Value = 126
$ws.Range("L29").Value = 619
$ws.Range("L30").Value = 57
$ws.Range("L33").Value = 541
$ws.Range("L36").Value = 157
$ws.Range("L37").Value = 411
$ws.Range("L40").Value = 33
$ws.Range("L41").Value = 51
$ws.Range("L43").Value = 86
$ws.Range("L44").Value = 83
$ws.Range("L48").Value = 161
$ws.Range("L51").Value = 145
$ws.Range("L53").Value = 127
$ws.Range("L54").Value = 244
$ws.Range("L60").Value = 68
$ws.Range("L62").Value = 5
$ws.Range("L63").Value = 36
$ws.Range("L65").Value = 222
$ws.Range("L67").Value = 411
$ws.Range("L72").Value = 52
$ws.Range("L76").Value = 168
$ws.Range("L77").Value = 69
$ws.Range("L79").Value = 299
$ws.Range("L83").Value = 266
$ws.Range("L84").Value = 114
$ws.Range("L85").Value = 584
$ws.Range("L86").Value = 87
$ws.Range("L89").Value = 160
$ws.Range("L90").Value = 110
$ws.Range("L91").Value = 164
$ws.Range("L93").Value = 61
$ws.Range("L95").Value = 164
$ws.Range("L99").Value = 192
$ws.Range("L101").Value = 11519

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 128
$ws.Range("L3").Value = 117
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 389

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 192

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 46
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 169
$ws.Range("L3").Value = 240
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 584

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 216
$ws.Range("L3").Value = 248
$ws.Range("L5").Value = 28
$ws.Range("L6").Value = 207
$ws.Range("L7").Value = 747

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 88
$ws.Range("L3").Value = 106
$ws.Range("L7").Value = 266

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 173
$ws.Range("L7").Value = 541

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 64
$ws.Range("L3").Value = 51
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 126
$ws.Range("L3").Value = 129
$ws.Range("L7").Value = 411

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 79
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 77
$ws.Range("L7").Value = 192

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 155
$ws.Range("L7").Value = 411

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 181
$ws.Range("L3").Value = 241
$ws.Range("L7").Value = 619

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 161

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 112
$ws.Range("L3").Value = 100
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 326

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 9
$ws.Range("L6").Value = 18

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 70
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 102
$ws.Range("L3").Value = 109
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 92
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 157

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L3").Value = 15
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 44
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 5
